# Mirror the first "Sprint Billing" cost table (rows 1-12) onto the
# second one (rows 14-24): fill in the Cost column (D) for each line item
# and add a "Total" row that sums the Hours and Cost columns.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sprint Billing")

# D17:D24 get the same per-row formatting as D4:D11 (plain number format
# for the line items, the shaded "total band" format on the last line).
$ws.Range("D4:D11").Copy()
$ws.Range("D17:D24").PasteSpecial(-4122)
$ws.Range("D17:D24").Value = 500

# New row 25 mirrors row 12 (the existing "Total" row): same shading,
# a label in column A, and SUM formulas in B and D.
$ws.Range("A12:D12").Copy()
$ws.Range("A25:D25").PasteSpecial(-4122)
$ws.Range("A25").Value = $ws.Range("A12").Value2
$ws.Range("B25").Formula = "=SUM(B16:B24)"
$ws.Range("D25").Formula = "=SUM(D17:D24)"

$excel.CutCopyMode = 0

# Restore the active selection to match the authored workbook.
$ws.Range("F19").Select()
